# Update the "Use Statement" sheet with the new MITRE release / distribution
# text (replacing the old copyright / public-release-number footer block).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Use Statement")

# A1 keeps the same text as before ("This is part of the Use Case Exercise 4
# materials") - set it explicitly anyway so the cell is not left dependent on
# whatever happened to be there already.
$ws.Range("A1").Value = "This is part of the Use Case Exercise 4 materials"

# A2: new distribution statement text.
$ws.Range("A2").Value = "Distribution Statement A.  Approved for public release: distribution is unlimited"

# A3: new MITRE public release case number text.
$ws.Range("A3").Value = "MITRE Public Release Case Number 24-2089"

# A4: rich text - a copyright glyph (Symbol font char 0xE3 renders as (c))
# followed by bold "2025 The MITRE Corporation." in Calibri.
$copyrightChar = [char]0x00E3
$text = "$copyrightChar 2025 The MITRE Corporation."
$ws.Range("A4").Value = $text

# Cell-level font acts as the format for the first (unstyled) character run -
# make it the bold Symbol font so the leading glyph renders as a copyright
# sign.
$ws.Range("A4").Font.Name = "Symbol"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Size = 11

# The remainder of the string (" 2025 The MITRE Corporation.") is explicitly
# formatted as bold Calibri 11.
$rest = $ws.Range("A4").Characters(2, $text.Length - 1)
$rest.Font.Name = "Calibri"
$rest.Font.Bold = $true
$rest.Font.Size = 11
